$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1. Update the "Boost.Format" performance measurements (used by both the
#    worksheet cells and the two bar charts that cache them).
# ---------------------------------------------------------------------------
$ws.Range("E6:J6").Value = @(3.1880000000000002, 3.7370000000000001, 2.8780000000000001, 3.2170000000000001, 2.6720000000000002, 2.0110000000000001)
$ws.Range("E39:J39").Value = @(6.0220000000000002, 6.8170000000000002, 5.7930000000000001, 5.7210000000000001, 4.3540000000000001, 4.0890000000000004)

# ---------------------------------------------------------------------------
# 2. Drop the stale external reference to the old palette-chooser workbook.
# ---------------------------------------------------------------------------
foreach ($link in @($wb.LinkSources(1))) {
    if ($link) {
        $wb.BreakLink($link, 1)
    }
}

# ---------------------------------------------------------------------------
# 3. Tweak the value axis of the first chart ("Format single double"):
#    major unit 2 -> 1, and drop the explicit minor unit (go back to auto).
# ---------------------------------------------------------------------------
$chart1 = $ws.ChartObjects(1).Chart
$valAx1 = $chart1.Axes(2)
$valAx1.MajorUnit = 1
$valAx1.MinorUnit = [System.Reflection.Missing]::Value

# ---------------------------------------------------------------------------
# 4. Tweak the second chart ("Format sequence of several items"): remove the
#    fixed maximum on the value axis so it goes back to automatic scaling.
# ---------------------------------------------------------------------------
$chart2 = $ws.ChartObjects(2).Chart
$valAx2 = $chart2.Axes(2)
$valAx2.MaximumScale = [System.Reflection.Missing]::Value

# ---------------------------------------------------------------------------
# 5. Update the view state of Sheet1 to match where the author was working.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E39:J39").Select()
